$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 4904349.5
$ws.Cells.Item(17, 10).Value = 4904349.5
$ws.Cells.Item(17, 12).Value = 14713048.5
$ws.Cells.Item(17, 14).Value = -14713384.5
$ws.Cells.Item(29, 8).Value = 100000
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 100000
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 300000
$ws.Cells.Item(29, 13).ClearContents()
$ws.Cells.Item(29, 14).Value = -300562
$ws.Cells.Item(69, 8).Value = 3651.625
$ws.Cells.Item(69, 10).Value = 3733.3333
$ws.Cells.Item(69, 12).Value = 11199.9999
$ws.Cells.Item(69, 14).Value = -12947.9999
$ws.Cells.Item(72, 8).Value = 3651.625
$ws.Cells.Item(72, 10).Value = 3733.3333
$ws.Cells.Item(72, 12).Value = 33599.9997
$ws.Cells.Item(72, 14).Value = -42335.9997
$ws.Cells.Item(113, 8).Value = 2666.6667
$ws.Cells.Item(113, 9).Value = 2615.3845
$ws.Cells.Item(113, 11).Value = 2615.3845
$ws.Cells.Item(113, 13).Value = 638.6154999999999
$ws.Cells.Item(116, 8).Value = 2882.0322
$ws.Cells.Item(116, 9).Value = 1529.2222
$ws.Cells.Item(116, 10).Value = 3435.4546
$ws.Cells.Item(116, 11).Value = 1529.2222
$ws.Cells.Item(116, 12).Value = 3435.4546
$ws.Cells.Item(116, 13).Value = 1912.7778
$ws.Cells.Item(116, 14).Value = -10319.4546
$ws.Cells.Item(118, 8).Value = 1492.4667
$ws.Cells.Item(118, 9).Value = 432
$ws.Cells.Item(118, 10).Value = 2022.7
$ws.Cells.Item(118, 11).Value = 1296
$ws.Cells.Item(118, 12).Value = 6068.1
$ws.Cells.Item(118, 13).Value = 361
$ws.Cells.Item(118, 14).Value = -9382.1
$ws.Cells.Item(125, 8).Value = 2606
$ws.Cells.Item(125, 9).Value = 4166.6665
$ws.Cells.Item(125, 10).Value = 1045.3334
$ws.Cells.Item(125, 11).Value = 37499.9985
$ws.Cells.Item(125, 12).Value = 9408.000599999999
$ws.Cells.Item(125, 13).Value = -35039.9985
$ws.Cells.Item(125, 14).Value = -14328.0006
$ws.Cells.Item(129, 8).Value = 4287.1113
$ws.Cells.Item(129, 10).Value = 5709.697
$ws.Cells.Item(129, 12).Value = 17129.091
$ws.Cells.Item(129, 14).Value = -27129.091
$ws.Cells.Item(132, 8).Value = 2190.6758
$ws.Cells.Item(132, 9).Value = 1794.2388
$ws.Cells.Item(132, 10).Value = 5985.143
$ws.Cells.Item(132, 11).Value = 5382.7164
$ws.Cells.Item(132, 12).Value = 17955.429
$ws.Cells.Item(132, 13).Value = -2852.7164
$ws.Cells.Item(132, 14).Value = -23015.429
$ws.Cells.Item(135, 8).Value = 1464.2084
$ws.Cells.Item(135, 9).Value = 1125.0488
$ws.Cells.Item(135, 10).Value = 3450.7144
$ws.Cells.Item(135, 11).Value = 10125.4392
$ws.Cells.Item(135, 12).Value = 31056.4296
$ws.Cells.Item(135, 13).Value = -7590.439200000001
$ws.Cells.Item(135, 14).Value = -36126.4296
$ws.Cells.Item(137, 8).Value = 8696642
$ws.Cells.Item(137, 9).Value = 969.5454999999999
$ws.Cells.Item(137, 10).Value = 16667675
$ws.Cells.Item(137, 11).Value = 2908.6365
$ws.Cells.Item(137, 12).Value = 50003025
$ws.Cells.Item(137, 13).Value = -358.6364999999996
$ws.Cells.Item(137, 14).Value = -50008125
$ws.Cells.Item(139, 8).Value = 49872.5
$ws.Cells.Item(139, 10).Value = 49872.5
$ws.Cells.Item(139, 12).Value = 49872.5
$ws.Cells.Item(139, 14).Value = -60152.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 9435185
$ws.Cells.Item(61, 9).Value = 12196298
$ws.Cells.Item(61, 10).Value = 1381.0834
$ws.Cells.Item(61, 11).Value = 12196298
$ws.Cells.Item(61, 12).Value = 1381.0834
$ws.Cells.Item(61, 13).Value = -12196086
$ws.Cells.Item(61, 14).Value = -1805.0834
$ws.Cells.Item(132, 8).Value = 5683790
$ws.Cells.Item(132, 9).Value = 8930355
$ws.Cells.Item(132, 10).Value = 2300.0625
$ws.Cells.Item(132, 11).Value = 26791065
$ws.Cells.Item(132, 12).Value = 6900.1875
$ws.Cells.Item(132, 13).Value = -26788535
$ws.Cells.Item(132, 14).Value = -11960.1875
$ws.Cells.Item(136, 8).Value = 9435185
$ws.Cells.Item(136, 9).Value = 12196298
$ws.Cells.Item(136, 10).Value = 1381.0834
$ws.Cells.Item(136, 11).Value = 36588894
$ws.Cells.Item(136, 12).Value = 4143.2502
$ws.Cells.Item(136, 13).Value = -36586344
$ws.Cells.Item(136, 14).Value = -9243.2502

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 4720
$ws.Cells.Item(105, 9).Value = 2900
$ws.Cells.Item(105, 10).Value = 5000
$ws.Cells.Item(105, 11).Value = 2900
$ws.Cells.Item(105, 12).Value = 5000
$ws.Cells.Item(105, 13).Value = -1153
$ws.Cells.Item(105, 14).Value = -8494

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 952.37933
$ws.Cells.Item(16, 9).Value = 956.3
$ws.Cells.Item(16, 10).Value = 943.6667
$ws.Cells.Item(16, 11).Value = 956.3
$ws.Cells.Item(16, 12).Value = 943.6667
$ws.Cells.Item(16, 13).Value = -669.3
$ws.Cells.Item(16, 14).Value = -1517.6667
$ws.Cells.Item(41, 8).Value = 26413
$ws.Cells.Item(41, 10).Value = 30016.25
$ws.Cells.Item(41, 12).Value = 30016.25
$ws.Cells.Item(41, 14).Value = -30872.25
$ws.Cells.Item(107, 8).Value = 874.53845
$ws.Cells.Item(107, 9).Value = 892.63635
$ws.Cells.Item(107, 10).Value = 775
$ws.Cells.Item(107, 11).Value = 892.63635
$ws.Cells.Item(107, 12).Value = 775
$ws.Cells.Item(107, 13).Value = 1027.36365
$ws.Cells.Item(107, 14).Value = -4615
$ws.Cells.Item(113, 8).Value = 952.37933
$ws.Cells.Item(113, 9).Value = 956.3
$ws.Cells.Item(113, 10).Value = 943.6667
$ws.Cells.Item(113, 11).Value = 956.3
$ws.Cells.Item(113, 12).Value = 943.6667
$ws.Cells.Item(113, 13).Value = 1213.7
$ws.Cells.Item(113, 14).Value = -5283.6667
$ws.Cells.Item(122, 8).Value = 1465.5769
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 14).Value = -10900
$ws.Cells.Item(132, 8).Value = 8476007
$ws.Cells.Item(132, 9).Value = 11112178
$ws.Cells.Item(132, 10).Value = 2601.8572
$ws.Cells.Item(132, 11).Value = 33336534
$ws.Cells.Item(132, 12).Value = 7805.571599999999
$ws.Cells.Item(132, 13).Value = -33334004
$ws.Cells.Item(132, 14).Value = -12865.5716
$ws.Cells.Item(140, 8).Value = 45722.715
$ws.Cells.Item(140, 10).Value = 45722.715
$ws.Cells.Item(140, 12).Value = 45722.715
$ws.Cells.Item(140, 14).Value = -56082.715

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 44
$ws.Cells.Item(2, 9).Value = 37.833332
$ws.Cells.Item(2, 10).Value = 53.25
$ws.Cells.Item(2, 11).Value = 37.833332
$ws.Cells.Item(2, 12).Value = 53.25
$ws.Cells.Item(2, 13).Value = 75.166668
$ws.Cells.Item(2, 14).Value = -279.25
$ws.Cells.Item(18, 8).Value = 13900
$ws.Cells.Item(18, 10).Value = 13900
$ws.Cells.Item(18, 12).Value = 13900
$ws.Cells.Item(18, 14).Value = -14486
$ws.Cells.Item(102, 8).Value = 5647.4614
$ws.Cells.Item(102, 9).Value = 5701.4165
$ws.Cells.Item(102, 11).Value = 5701.4165
$ws.Cells.Item(102, 13).Value = -4079.4165
$ws.Cells.Item(107, 8).Value = 3445.6667
$ws.Cells.Item(107, 9).Value = 3445.6667
$ws.Cells.Item(107, 11).Value = 3445.6667
$ws.Cells.Item(107, 13).Value = -1525.6667
$ws.Cells.Item(113, 8).Value = 144537.42
$ws.Cells.Item(113, 9).Value = 144537.42
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 144537.42
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -142367.42
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 7410595
$ws.Cells.Item(122, 9).Value = 16668841
$ws.Cells.Item(122, 10).Value = 3998
$ws.Cells.Item(122, 11).Value = 50006523
$ws.Cells.Item(122, 12).Value = 11994
$ws.Cells.Item(122, 13).Value = -50004073
$ws.Cells.Item(122, 14).Value = -16894
$ws.Cells.Item(126, 8).Value = 4332.914
$ws.Cells.Item(126, 9).Value = 3096.5293
$ws.Cells.Item(126, 10).Value = 5500.6113
$ws.Cells.Item(126, 11).Value = 9289.5879
$ws.Cells.Item(126, 12).Value = 16501.8339
$ws.Cells.Item(126, 13).Value = -6819.5879
$ws.Cells.Item(126, 14).Value = -21441.8339
$ws.Cells.Item(132, 8).Value = 3146.9812
$ws.Cells.Item(132, 9).Value = 2256.95
$ws.Cells.Item(132, 10).Value = 5885.5386
$ws.Cells.Item(132, 11).Value = 6770.849999999999
$ws.Cells.Item(132, 12).Value = 17656.6158
$ws.Cells.Item(132, 13).Value = -4240.849999999999
$ws.Cells.Item(132, 14).Value = -22716.6158

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5156.288
$ws.Cells.Item(7, 9).Value = 5197.8887
$ws.Cells.Item(7, 10).Value = 5091.174
$ws.Cells.Item(7, 11).Value = 5197.8887
$ws.Cells.Item(7, 12).Value = 5091.174
$ws.Cells.Item(7, 13).Value = -5085.8887
$ws.Cells.Item(7, 14).Value = -5315.174
$ws.Cells.Item(61, 8).Value = 1429.4706
$ws.Cells.Item(61, 9).Value = 1311.9286
$ws.Cells.Item(61, 10).Value = 1978
$ws.Cells.Item(61, 11).Value = 1311.9286
$ws.Cells.Item(61, 12).Value = 1978
$ws.Cells.Item(61, 13).Value = -1109.9286
$ws.Cells.Item(61, 14).Value = -2382
$ws.Cells.Item(113, 8).Value = 1429.4706
$ws.Cells.Item(113, 9).Value = 1311.9286
$ws.Cells.Item(113, 10).Value = 1978
$ws.Cells.Item(113, 11).Value = 1311.9286
$ws.Cells.Item(113, 12).Value = 1978
$ws.Cells.Item(113, 13).Value = 858.0714
$ws.Cells.Item(113, 14).Value = -6318
$ws.Cells.Item(122, 8).Value = 4058.121
$ws.Cells.Item(122, 9).Value = 4201.0454
$ws.Cells.Item(122, 11).Value = 12603.1362
$ws.Cells.Item(122, 13).Value = -10153.1362
$ws.Cells.Item(126, 8).Value = 5156.288
$ws.Cells.Item(126, 9).Value = 5197.8887
$ws.Cells.Item(126, 10).Value = 5091.174
$ws.Cells.Item(126, 11).Value = 15593.6661
$ws.Cells.Item(126, 12).Value = 15273.522
$ws.Cells.Item(126, 13).Value = -13123.6661
$ws.Cells.Item(126, 14).Value = -20213.522
$ws.Cells.Item(132, 8).Value = 8587.056
$ws.Cells.Item(132, 9).Value = 4410.55
$ws.Cells.Item(132, 10).Value = 13807.6875
$ws.Cells.Item(132, 11).Value = 13231.65
$ws.Cells.Item(132, 12).Value = 41423.0625
$ws.Cells.Item(132, 13).Value = -10701.65
$ws.Cells.Item(132, 14).Value = -46483.0625
$ws.Cells.Item(136, 8).Value = 16676198
$ws.Cells.Item(136, 9).Value = 27780774
$ws.Cells.Item(136, 10).Value = 19334.166
$ws.Cells.Item(136, 11).Value = 83342322
$ws.Cells.Item(136, 12).Value = 58002.49800000001
$ws.Cells.Item(136, 13).Value = -83339772
$ws.Cells.Item(136, 14).Value = -63102.49800000001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 8058.68
$ws.Cells.Item(62, 9).Value = 4850.933
$ws.Cells.Item(62, 10).Value = 12870.3
$ws.Cells.Item(62, 11).Value = 4850.933
$ws.Cells.Item(62, 12).Value = 12870.3
$ws.Cells.Item(62, 13).Value = -4226.933
$ws.Cells.Item(62, 14).Value = -14118.3
$ws.Cells.Item(65, 8).Value = 8058.68
$ws.Cells.Item(65, 9).Value = 4850.933
$ws.Cells.Item(65, 10).Value = 12870.3
$ws.Cells.Item(65, 11).Value = 24254.665
$ws.Cells.Item(65, 12).Value = 64351.5
$ws.Cells.Item(65, 13).Value = -21134.665
$ws.Cells.Item(65, 14).Value = -70591.5
$ws.Cells.Item(122, 8).Value = 5484.6665
$ws.Cells.Item(122, 9).Value = 5484.6665
$ws.Cells.Item(122, 11).Value = 16453.9995
$ws.Cells.Item(122, 13).Value = -14003.9995
$ws.Cells.Item(126, 8).Value = 9638.5
$ws.Cells.Item(126, 9).Value = 3952
$ws.Cells.Item(126, 10).Value = 15325
$ws.Cells.Item(126, 11).Value = 11856
$ws.Cells.Item(126, 12).Value = 45975
$ws.Cells.Item(126, 13).Value = -9386
$ws.Cells.Item(126, 14).Value = -50915
$ws.Cells.Item(132, 8).Value = 1720.1698
$ws.Cells.Item(132, 9).Value = 1459.1904
$ws.Cells.Item(132, 10).Value = 2716.6365
$ws.Cells.Item(132, 11).Value = 4377.5712
$ws.Cells.Item(132, 12).Value = 8149.9095
$ws.Cells.Item(132, 13).Value = -1847.5712
$ws.Cells.Item(132, 14).Value = -13209.9095
